$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RACP")

# Update the formula in B2 to round the result to 0 decimal places
$ws.Range("B2").Formula = "=ROUND(100/About!A11,0)"
